# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1068
    "F5"  = 3061
    "F7"  = 2280
    "F11" = 1089
    "F13" = 45
    "F15" = 227
    "F16" = 273
    "F18" = 11
    "F19" = 8
    "F20" = 99
    "F21" = 49
    "F22" = 74
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
